$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data rows (22-31) need to shift down by two rows to make
# room for two new "Dina" price records inserted at the top of the block
# (new rows 22 and 23). Inserting whole rows preserves the formatting
# (date style, etc.) of the rows being pushed down.
$ws.Rows("22:23").Insert()

# --- New row 22: Dina / Especial, Región Metropolitana, $/caja 16 kilos ---
$ws.Cells.Item(22, 1).Value = 2
$ws.Cells.Item(22, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(22, 3).Value = "Coquimbo"
$ws.Cells.Item(22, 4).Value2 = 44917
$ws.Cells.Item(22, 5).Value = 4
$ws.Cells.Item(22, 6).Value = "Fruta"
$ws.Cells.Item(22, 7).Value = 100103
$ws.Cells.Item(22, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(22, 9).Value = 100103003
$ws.Cells.Item(22, 10).Value = "Damasco"
$ws.Cells.Item(22, 11).Value = "Dina"
$ws.Cells.Item(22, 12).Value = "Especial"
$ws.Cells.Item(22, 13).Value = 100
$ws.Cells.Item(22, 14).Value = 23000
$ws.Cells.Item(22, 15).Value = 24000
$ws.Cells.Item(22, 16).Value = 23500
$ws.Cells.Item(22, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(22, 18).Value = "Región Metropolitana"
$ws.Cells.Item(22, 19).Value = 1469
$ws.Cells.Item(22, 20).Value = 16

# --- New row 23: Dina / Primera, Región Metropolitana, $/caja 16 kilos ---
$ws.Cells.Item(23, 1).Value = 2
$ws.Cells.Item(23, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(23, 3).Value = "Coquimbo"
$ws.Cells.Item(23, 4).Value2 = 44917
$ws.Cells.Item(23, 5).Value = 4
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100103
$ws.Cells.Item(23, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(23, 9).Value = 100103003
$ws.Cells.Item(23, 10).Value = "Damasco"
$ws.Cells.Item(23, 11).Value = "Dina"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 100
$ws.Cells.Item(23, 14).Value = 20000
$ws.Cells.Item(23, 15).Value = 21000
$ws.Cells.Item(23, 16).Value = 20500
$ws.Cells.Item(23, 17).Value = "$/caja 16 kilos"
$ws.Cells.Item(23, 18).Value = "Región Metropolitana"
$ws.Cells.Item(23, 19).Value = 1281
$ws.Cells.Item(23, 20).Value = 16
